$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Update the receipt/record id and the id-card number for the second record
$ws.Range("A2").Value = 3023
$ws.Range("E2").Value = 46200608023

# Reflect the last-used cell selection on the Data sheet
$ws.Range("I9").Select()
